# DCF_Form_9.xlsx - "guidelines in template import excel"
#
# Adds three new guideline rows (Type of Training/Activity, Name of
# Partner/Organization, Counterpart Amount) with their accompanying
# "how to fill this in" notes to the Guidelines sheet, widens the two
# columns on that sheet so the new text is readable, and moves the
# active selection the way the author left it after editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Guidelines")

# ---------------------------------------------------------------------
# 1. New guidance rows (3, 4, 5)
# ---------------------------------------------------------------------

# Row 3 - "Type of Training/Activity" guideline. Column A reuses the
# green "field name" look of A2, column B reuses the tan "guideline
# text" look of B2 but additionally needs word-wrap since the note is
# long.
$ws.Range("A2").Copy()
$ws.Range("A3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A3").Value = "Type of Training/Activity"

$ws.Range("B2").Copy()
$ws.Range("B3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("B3").Value = "Select one (Capbuild,Meetings,Policy issuances,Others) if others please specify"
$ws.Range("B3").WrapText = $true

# Row 4 - "Name of Partner/Organization"
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "Name of Partner/Organization"

$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B4").Value = "if multiple data add comma (ex. sample text1, sample text2, sample text3)"

# Row 5 - "Counterpart Amount"
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "Counterpart Amount"

$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$ws.Range("B5").Value = "if multiple data add comma (ex. sample1, sample2, sample3)"

# The green field-name cells (A3:A5) use a plain black font rather than
# the automatic theme color the copied format carries.
$ws.Range("A3:A5").Font.Color = 0

# Row heights to fit the (now taller / wrapped) guideline rows.
$ws.Rows.Item(3).RowHeight = 30.75
$ws.Rows.Item(4).RowHeight = 15.75
$ws.Rows.Item(5).RowHeight = 15.75

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2. Column widths - widened to fit the new guideline text
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 34
$ws.Columns.Item(2).ColumnWidth = 66.5

# ---------------------------------------------------------------------
# 3. Leave the selection where the author left it
# ---------------------------------------------------------------------
$ws.Range("B13").Select()
